$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header B1 text and add new column G header
$ws.Range("B1").Value = "Core Profile Title = id"
$ws.Range("G1").Value = "code"
$ws.Range("G1").Font.Bold = $true

# Update column B values (Core Profile Title = id) and add column G (code) for each data row
$ws.Range("B2").Value = "AllergyIntolerance"
$ws.Range("G2").Value = "core-allergyintolerance"
$ws.Range("B3").Value = "CarePlan"
$ws.Range("G3").Value = "core-careplan"
$ws.Range("B4").Value = "CareTeam"
$ws.Range("G4").Value = "core-careteam"
$ws.Range("B5").Value = "Condition"
$ws.Range("G5").Value = "core-condition"
$ws.Range("B6").Value = "Conformance"
$ws.Range("G6").Value = "core-conformance"
$ws.Range("B7").Value = "Device"
$ws.Range("G7").Value = "core-device"
$ws.Range("B8").Value = "DiagnosticReport-Results"
$ws.Range("G8").Value = "core-diagnosticreport"
$ws.Range("B9").Value = "DocumentReference"
$ws.Range("G9").Value = "core-documentreference"
$ws.Range("B10").Value = "Goals"
$ws.Range("G10").Value = "core-goals"
$ws.Range("B11").Value = "Immunization"
$ws.Range("G11").Value = "core-immunization"
$ws.Range("B12").Value = "Location"
$ws.Range("G12").Value = "core-location"
$ws.Range("B13").Value = "Medication"
$ws.Range("G13").Value = "core-medication"
$ws.Range("B14").Value = "Medicationadministration"
$ws.Range("G14").Value = "core-medicationadministration"
$ws.Range("B15").Value = "Medicationdispense"
$ws.Range("G15").Value = "core-medicationdispense"
$ws.Range("B16").Value = "Medicationorder"
$ws.Range("G16").Value = "core-medicationorder"
$ws.Range("B17").Value = "Medicationstatement"
$ws.Range("G17").Value = "core-medicationstatement"
$ws.Range("B18").Value = "Observation-Results"
$ws.Range("G18").Value = "core-observation-results"
$ws.Range("B19").Value = "Observation-Resultsv2"
$ws.Range("G19").Value = "core-observation-resultsv2"
$ws.Range("B20").Value = "Observation-Smokingstatus"
$ws.Range("G20").Value = "core-observation-smokingstatus"
$ws.Range("B21").Value = "Observation-Vitalsigns"
$ws.Range("G21").Value = "core-observation-vitalsigns"
$ws.Range("B22").Value = "Organization"
$ws.Range("G22").Value = "core-organization"
$ws.Range("B23").Value = "Patient"
$ws.Range("G23").Value = "core-patient"
$ws.Range("B24").Value = "Practitioner"
$ws.Range("G24").Value = "core-practitioner"
$ws.Range("B25").Value = "Procedure"
$ws.Range("G25").Value = "core-procedure"

# Re-point the active selection to B1 (was B2)
$ws.Range("B1").Select()
